# Corrected Calibration and Ingest Sheets for Coastal Gliders
# - FLORT cal value "CC_angular_resolution" (row 9) -> 1.076
# - FLORT cal value "CC_scattering_angle" (row 7) -> 124
# Also reflects that the workbook was last left with the
# "Asset_Cal_Info" tab active/selected (cell F33) rather than "Moorings".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

# Update the FLORT calibration coefficients.
$ws.Range("F7").Value = 124
$ws.Range("F9").Value = 1.076

# Make Asset_Cal_Info the active sheet/tab with F33 selected, matching
# the saved view state of the workbook after editing.
$ws.Activate()
$ws.Range("F33").Select()
